$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to Text first so number-like strings (e.g. "1.002")
# are stored as text (matching the source t="inlineStr" cells) rather than being
# auto-parsed into numbers. Revert the format afterwards so cell styling (s attr)
# is left exactly as it was.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.115.73"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.883.46"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "322.40"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.4691"
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("D8").Value = "0.4016"
$ws.Range("E8").Value = "  -2.33%  "
$ws.Range("D9").Value = "47.38"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "0.07997"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").Value = "0.9909"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "22.32"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "1.890.94"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "5.844"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "7.013"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "88.67"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "0.06629"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "0.00001026"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "17.39"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "29.123.93"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "5.476"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "11.60"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "2.182"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "2.121.36"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "155.26"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "19.57"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "6.070"
$ws.Range("E29").Value = "  +8.43%  "
$ws.Range("D30").Value = "2.064"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").Value = "117.13"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "1.036"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").Value = "0.09423"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").Value = "3.542"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").Value = "1.376"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").Value = "5.334"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "0.06061"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "0.02224"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").Value = "1.173"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "8.011"
$ws.Range("E40").Value = "  -4.20%  "
$ws.Range("D41").Value = "0.5792"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "0.1819"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "2.462"
$ws.Range("E43").Value = "  +7.70%  "
$ws.Range("D44").Value = "9.998"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "1.271"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").Value = "0.07658"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "0.5457"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "1.895"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "113.36"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "43.93"
$ws.Range("E51").Value = "  +0.61%  "

$priceRange.Style = "Normal"
